$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Replace the curly "smart" quotes around the Tier labels (column H)
#    with plain straight quotes, for every row that shows them.
# ------------------------------------------------------------------
$ws.Range("H9").Value  = """Tier1"""
$ws.Range("H10").Value = """Tier1"""
$ws.Range("H11").Value = """Tier2"""
$ws.Range("H12").Value = """Tier2"""
$ws.Range("H13").Value = """Tier3"""
$ws.Range("H14").Value = """Tier3"""
$ws.Range("H15").Value = """Tier4"""
$ws.Range("H16").Value = """Tier4"""
$ws.Range("H17").Value = """Tier5"""
$ws.Range("H18").Value = """Tier5"""

# ------------------------------------------------------------------
# 2) Re-apply the cell formatting (font + fill + centered alignment)
#    on H10:H18 -- this is what produced the extra font/style entry
#    in the saved workbook.
# ------------------------------------------------------------------
$fmtRange = $ws.Range("H10:H18")
$fmtRange.Font.Name = "Calibri"
$fmtRange.Font.Size = 11
$fmtRange.Font.Color = 0
$fmtRange.Font.Bold = $false
$fmtRange.Interior.Pattern = 1
$fmtRange.Interior.Color = $ws.Range("H9").Interior.Color
$fmtRange.HorizontalAlignment = -4108
$fmtRange.VerticalAlignment = -4108

# ------------------------------------------------------------------
# 3) Move the active selection from I9 to A19 (last thing the author
#    did before saving).
# ------------------------------------------------------------------
$ws.Range("A19").Select()
